# Update "countries & provincias Spain" COVID data sheet.
# The sheet (single worksheet "Pais") lists countries sorted descending by
# "Casos totales" (column B). This update refreshes several countries'
# daily figures; because some countries' totals now overtake their
# neighbours, those rows swap positions (country name + full data row
# moves as a unit) to keep the list sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 13:14"

# --- helper: write one full data row (B:H) ---
# Row 6: India - values refreshed, stays in place
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1390429
$ws.Range("C6").Value = 4935
$ws.Range("D6").Value = 889187
$ws.Range("E6").Value = 469091
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 55
$ws.Range("H6").Value = 32151

# Rows 42-43: Emiratos Arabes Unidos overtakes Panama
$ws.Range("A42").Value = "Emiratos Arabes Unidos"
$ws.Range("B42").Value = 58913
$ws.Range("C42").Value = 351
$ws.Range("D42").Value = 52182
$ws.Range("E42").Value = 6387
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 344

$ws.Range("A43").Value = "Panama"
$ws.Range("B43").Value = 58864
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 33428
$ws.Range("E43").Value = 24161
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 1275

# Rows 47-48: Rumania overtakes Guatemala
$ws.Range("A47").Value = "Rumania"
$ws.Range("B47").Value = 44798
$ws.Range("C47").Value = 1120
$ws.Range("D47").Value = 25643
$ws.Range("E47").Value = 16968
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 22
$ws.Range("H47").Value = 2187

$ws.Range("A48").Value = "Guatemala"
$ws.Range("B48").Value = 44492
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 31045
$ws.Range("E48").Value = 11748
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 1699

# Rows 85-87: Madagascar overtakes Noruega and Malasia
$ws.Range("A85").Value = "Madagascar"
$ws.Range("B85").Value = 9295
$ws.Range("C85").Value = 429
$ws.Range("D85").Value = 5579
$ws.Range("E85").Value = 3631
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 7
$ws.Range("H85").Value = 85

$ws.Range("A86").Value = "Noruega"
$ws.Range("B86").Value = 9111
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 8674
$ws.Range("E86").Value = 182
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 255

$ws.Range("A87").Value = "Malasia"
$ws.Range("B87").Value = 8897
$ws.Range("C87").Value = 13
$ws.Range("D87").Value = 8600
$ws.Range("E87").Value = 173
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 124

# Row 114: Sri Lanka - values refreshed, stays in place
$ws.Range("A114").Value = "Sri Lanka"
$ws.Range("B114").Value = 2772
$ws.Range("C114").Value = 2
$ws.Range("D114").Value = 2106
$ws.Range("E114").Value = 655
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 11

# Rows 115-116: Montenegro overtakes Hong Kong
$ws.Range("A115").Value = "Montenegro"
$ws.Range("B115").Value = 2747
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 664
$ws.Range("E115").Value = 2040
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 43

$ws.Range("A116").Value = "Hong Kong"
$ws.Range("B116").Value = 2634
$ws.Range("C116").Value = 128
$ws.Range("D116").Value = 1495
$ws.Range("E116").Value = 1121
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 18

# Rows 130-132: Benin overtakes Sierra Leona and Ruanda
$ws.Range("A130").Value = "Benin"
$ws.Range("B130").Value = 1770
$ws.Range("C130").Value = 76
$ws.Range("D130").Value = 1036
$ws.Range("E130").Value = 699
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 35

$ws.Range("A131").Value = "Sierra Leona"
$ws.Range("B131").Value = 1768
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 1297
$ws.Range("E131").Value = 405
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 66

$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 1752
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 907
$ws.Range("E132").Value = 840
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 5

# Rows 155-157: Malta overtakes San Marino and Botsuana
$ws.Range("A155").Value = "Malta"
$ws.Range("B155").Value = 700
$ws.Range("C155").Value = 14
$ws.Range("D155").Value = 665
$ws.Range("E155").Value = 26
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 9

$ws.Range("A156").Value = "San Marino"
$ws.Range("B156").Value = 699
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 657
$ws.Range("E156").Value = 0
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 42

$ws.Range("A157").Value = "Botsuana"
$ws.Range("B157").Value = 686
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 52
$ws.Range("E157").Value = 633
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 1

# Rows 162-163: Vietnam overtakes Lesoto
$ws.Range("A162").Value = "Vietnam"
$ws.Range("B162").Value = 420
$ws.Range("C162").Value = 3
$ws.Range("D162").Value = 365
$ws.Range("E162").Value = 55
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 0

$ws.Range("A163").Value = "Lesoto"
$ws.Range("B163").Value = 419
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 69
$ws.Range("E163").Value = 341
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 9

# Rows 187-189: Islas Turcas y Caicos overtakes San Martin (Parte Holandesa) and Butan
$ws.Range("A187").Value = "Islas Turcas y Caicos"
$ws.Range("B187").Value = 97
$ws.Range("C187").Value = 5
$ws.Range("D187").Value = 35
$ws.Range("E187").Value = 60
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 2

$ws.Range("A188").Value = "San Martin (Parte Holandesa)"
$ws.Range("B188").Value = 93
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 63
$ws.Range("E188").Value = 15
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 15

$ws.Range("A189").Value = "Butan"
$ws.Range("B189").Value = 93
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 85
$ws.Range("E189").Value = 8
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 0

# Rows 210-211: Islas Malvinas <-> Groenlandia swap position (figures tied, unchanged)
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Groenlandia"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
